# Atualizado por script em 13-11-2023 22:16
#
# The source site re-scraped the fixture list and a handful of matches that
# were previously written to the wrong row (two games kicking off on the
# same day got their row order swapped) needed to be corrected, plus one
# brand-new match (played 13/11/2023) had to be appended at the end of the
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- helper: swap the match-data columns (F..V) between two rows --------
# Columns A..E (Indice, pais, torneio, temporada, data_partida) are left
# untouched - only the home/away teams, odds, timestamps and match URL
# (columns F through V) need to trade places between the two rows.
function Swap-MatchRow($row1, $row2) {
    for ($col = 6; $col -le 22; $col++) {
        $c1 = $ws.Cells.Item($row1, $col)
        $c2 = $ws.Cells.Item($row2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value2 = $v2
        $c2.Value2 = $v1
    }
}

Swap-MatchRow 80 81
Swap-MatchRow 82 83
Swap-MatchRow 84 85
Swap-MatchRow 102 103
Swap-MatchRow 110 111

# --- append the new match row (row 182) ----------------------------------
$newRow = 182

# Copy the A/E number formats (bold+border style on the index column, the
# datetime display format on the match-date column) down from the last
# existing row before filling in the values.
$ws.Range("A181").Copy()
$ws.Range("A" + $newRow).PasteSpecial(-4122)
$ws.Range("E181").Copy()
$ws.Range("E" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item($newRow, 1).Value2 = 181
$ws.Cells.Item($newRow, 2).Value2 = "argentina"
$ws.Cells.Item($newRow, 3).Value2 = "copa-de-la-liga-profesional"
$ws.Cells.Item($newRow, 4).Value2 = "2023"
$ws.Cells.Item($newRow, 5).Value2 = 45243.95833333334
$ws.Cells.Item($newRow, 6).Value2 = "Sarmiento Junin"
$ws.Cells.Item($newRow, 7).Value2 = 0
$ws.Cells.Item($newRow, 8).Value2 = "Godoy Cruz"
$ws.Cells.Item($newRow, 9).Value2 = 0
$ws.Cells.Item($newRow, 10).Value2 = 3.01
$ws.Cells.Item($newRow, 11).Value2 = "07/11/2023 05:11"
$ws.Cells.Item($newRow, 12).Value2 = 3.27
$ws.Cells.Item($newRow, 13).Value2 = "13/11/2023 22:57"
$ws.Cells.Item($newRow, 14).Value2 = 3.02
$ws.Cells.Item($newRow, 15).Value2 = "07/11/2023 05:11"
$ws.Cells.Item($newRow, 16).Value2 = 2.85
$ws.Cells.Item($newRow, 17).Value2 = "13/11/2023 22:59"
$ws.Cells.Item($newRow, 18).Value2 = 2.71
$ws.Cells.Item($newRow, 19).Value2 = "07/11/2023 05:11"
$ws.Cells.Item($newRow, 20).Value2 = 2.66
$ws.Cells.Item($newRow, 21).Value2 = "13/11/2023 22:57"
$ws.Cells.Item($newRow, 22).Value2 = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/sarmiento-junin-godoy-cruz/vmmcZCe6/"
